$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tasks 04-08 to 04-21")
$ws.Rows.Item(7).Delete()
$ws.Rows.Item(10).Delete()
